# Add spamassassin setup codes
#
# - 연람희's "지난주 업무" cell (merged B12:C12) gets the new status note
#   about wiring Spamassassin into sendmail and scoping how to hook a
#   deep-learning model into it.
# - The "금주 업무" week-start date in B5 moves forward one day
#   (2021-09-11 -> 2021-09-12, serial 44450 -> 44451).
# - The sheet was left scrolled to row 5 with B12:C12 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report text for 연람희 in the "지난주 업무" column (merged B12:C12)
$ws.Range("B12").Value = " - Spamassassin sendmail 연동. Spamassassin에 딥러닝을 어떻게 연결할 수 있는지 분석"

# Bump the week date by one day
$ws.Range("B5").Value = 44451

# Match the saved view state: scrolled so row 5 is at the top, with the
# B12:C12 merged cell selected
$window = $excel.ActiveWindow
$window.ScrollRow = 5
$window.ScrollColumn = 1
$ws.Range("B12:C12").Select()
